$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete column B, shifting columns C:J (and everything to the right) one
# column to the left. Excel automatically adjusts cell references inside
# formulas when a column is deleted.
$ws.Columns("B:B").Delete()

# Update the active selection to match the post-edit state captured in the
# saved file (selection moved to F12 after the edit).
$ws.Range("F12").Select()
